$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the [Quantitation] and [Flags] sections to the manifest template
# below the existing [Publication] block (rows 1-29).
# NOTE: cell writes are intentionally ordered to reproduce the exact
# shared-string table ordering of the target workbook.

# [Quantitation] section
$ws.Range("A31").Value = "[Quantitation]"
$ws.Range("A31").Font.Bold = $true

$ws.Range("A32").Value = "#Channel"
$ws.Range("B32").Value = "#Identifier"

$ws.Range("A33").Value = "M"
$ws.Range("A34").Value = "L"
$ws.Range("B34").Value = "wt"
$ws.Range("B33").Value = "ko"

$ws.Range("A35").Value = "H"
$ws.Range("B35").Value = "sc"

# [Flags] section
$ws.Range("A37").Value = "[Flags]"
$ws.Range("A37").Font.Bold = $true

$ws.Range("B38").Value = "#Value"
$ws.Range("A38").Value = "#Flag"

$ws.Range("A39").Value = "hcd-alias"
$ws.Range("B39").Value = "PQD"

$ws.Range("A40").Value = "etd-alias"
$ws.Range("B40").Value = "CID"

$ws.Range("A42").Value = "feature_enable_truncated_nglyco"
$ws.Range("A41").Value = "hcd-dont-use-masspeaks"

$ws.Range("B41").Value = $true
$ws.Range("B42").Value = $true

$ws.Range("A41").Select()
